$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-13 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-14 Wednesday", 2)

$d.Content.Find.Execute("548÷4=137, 0", $true, $false, $false, $false, $false, $true, 1, $false, "915÷3=305, 0", 2)
$d.Content.Find.Execute("851÷7=121, 4", $true, $false, $false, $false, $false, $true, 1, $false, "891÷2=445, 1", 2)
$d.Content.Find.Execute("792÷5=158, 2", $true, $false, $false, $false, $false, $true, 1, $false, "740÷2=370, 0", 2)
$d.Content.Find.Execute("882÷2=441, 0", $true, $false, $false, $false, $false, $true, 1, $false, "829÷9=92, 1", 2)
$d.Content.Find.Execute("230÷6=38, 2", $true, $false, $false, $false, $false, $true, 1, $false, "293÷3=97, 2", 2)

$d.Content.Find.Execute("184÷3=61, 1", $true, $false, $false, $false, $false, $true, 1, $false, "332÷8=41, 4", 2)
$d.Content.Find.Execute("593÷7=84, 5", $true, $false, $false, $false, $false, $true, 1, $false, "546÷2=273, 0", 2)
$d.Content.Find.Execute("190÷8=23, 6", $true, $false, $false, $false, $false, $true, 1, $false, "290÷7=41, 3", 2)
$d.Content.Find.Execute("808÷8=101, 0", $true, $false, $false, $false, $false, $true, 1, $false, "857÷2=428, 1", 2)
$d.Content.Find.Execute("425÷4=106, 1", $true, $false, $false, $false, $false, $true, 1, $false, "854÷3=284, 2", 2)

$d.Content.Find.Execute("622÷7=88, 6", $true, $false, $false, $false, $false, $true, 1, $false, "376÷8=47, 0", 2)
$d.Content.Find.Execute("862÷6=143, 4", $true, $false, $false, $false, $false, $true, 1, $false, "856÷8=107, 0", 2)
$d.Content.Find.Execute("278÷4=69, 2", $true, $false, $false, $false, $false, $true, 1, $false, "576÷9=64, 0", 2)
$d.Content.Find.Execute("429÷8=53, 5", $true, $false, $false, $false, $false, $true, 1, $false, "986÷9=109, 5", 2)
$d.Content.Find.Execute("720÷4=180, 0", $true, $false, $false, $false, $false, $true, 1, $false, "660÷8=82, 4", 2)

$d.Content.Find.Execute("410÷6=68, 2", $true, $false, $false, $false, $false, $true, 1, $false, "296÷6=49, 2", 2)
$d.Content.Find.Execute("183÷4=45, 3", $true, $false, $false, $false, $false, $true, 1, $false, "868÷2=434, 0", 2)
$d.Content.Find.Execute("927÷4=231, 3", $true, $false, $false, $false, $false, $true, 1, $false, "262÷2=131, 0", 2)
$d.Content.Find.Execute("467÷2=233, 1", $true, $false, $false, $false, $false, $true, 1, $false, "850÷3=283, 1", 2)
$d.Content.Find.Execute("668÷3=222, 2", $true, $false, $false, $false, $false, $true, 1, $false, "822÷7=117, 3", 2)

$d.Content.Find.Execute("656÷2=328, 0", $true, $false, $false, $false, $false, $true, 1, $false, "966÷9=107, 3", 2)
$d.Content.Find.Execute("521÷8=65, 1", $true, $false, $false, $false, $false, $true, 1, $false, "929÷3=309, 2", 2)
$d.Content.Find.Execute("993÷5=198, 3", $true, $false, $false, $false, $false, $true, 1, $false, "807÷7=115, 2", 2)
$d.Content.Find.Execute("992÷7=141, 5", $true, $false, $false, $false, $false, $true, 1, $false, "430÷8=53, 6", 2)
$d.Content.Find.Execute("886÷2=443, 0", $true, $false, $false, $false, $false, $true, 1, $false, "994÷5=198, 4", 2)
